$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# Original: cols A:C = 36.42578125, cols D:H = 8.7109375, cols I:XFD = 9.140625 (default)
# Target:   cols A:C = 36.140625,  cols D:T = 8.5703125,  cols U:XFD = 9.140625 (default)
# (COM ColumnWidth is quantized to a pixel grid, so we use the closest
#  reachable values.)
$ws.Columns("A:C").ColumnWidth = 35.33
$ws.Columns("D:T").ColumnWidth = 7.67

# --- Row heights ---------------------------------------------------------
# Rows 2-4 pick up an explicit 13.5pt custom height.
$ws.Rows(2).RowHeight = 13.5
$ws.Rows(3).RowHeight = 13.5
$ws.Rows(4).RowHeight = 13.5

# --- New column T (year 2023) --------------------------------------------
# Row 3 is an empty, bottom-bordered filler cell - just clone S3's format.
$ws.Range("S3").Copy()
$ws.Range("T3").PasteSpecial(-4122)

# Row 4 header year.
$ws.Range("T4").Value = 2023
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)

# Row 5 data value.
$ws.Range("T5").Value = 1.4
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)

# Row 6 data value.
$ws.Range("T6").Value = 8.1999999999999993
$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)

# --- Selection -------------------------------------------------------------
# The old workbook had S3 selected; the new one clears that stale selection.
$ws.Range("A1").Select()
